$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 17.41485066666667
$ws.Cells.Item(2, 8).Value = 52.24455200000001
$ws.Cells.Item(2, 9).Value = 0.1047285618770465
$ws.Cells.Item(2, 10).Value = 0.1047285618770465
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.06447966666666667
$ws.Cells.Item(2, 14).Value = 0.193439
$ws.Cells.Item(2, 15).Value = 0.001101138907643723
$ws.Cells.Item(2, 16).Value = 0.001101138907643722
$ws.Cells.Item(2, 17).Value = 1.122903766036445
$ws.Cells.Item(2, 18).Value = 10.106133894328
$ws.Cells.Item(2, 19).Value = 0.000115320694224389
$ws.Cells.Item(2, 20).Value = 0.000115320694224389
$ws.Cells.Item(3, 7).Value = 17.41485066666667
$ws.Cells.Item(3, 8).Value = 52.24455200000001
$ws.Cells.Item(3, 9).Value = 0.1047285618770465
$ws.Cells.Item(3, 10).Value = 0.1047285618770465
$ws.Cells.Item(3, 15).Value = 0.00657695954769643
$ws.Cells.Item(3, 16).Value = 0.006576959547696431
$ws.Cells.Item(3, 17).Value = 6.706958217452446
$ws.Cells.Item(3, 18).Value = 60.36262395707201
$ws.Cells.Item(3, 19).Value = 0.0006887955149537576
$ws.Cells.Item(3, 20).Value = 0.0006887955149537577
$ws.Cells.Item(4, 7).Value = 17.41485066666667
$ws.Cells.Item(4, 8).Value = 52.24455200000001
$ws.Cells.Item(4, 9).Value = 0.1047285618770465
$ws.Cells.Item(4, 10).Value = 0.1047285618770465
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.008175
$ws.Cells.Item(4, 14).Value = 0.024525
$ws.Cells.Item(4, 15).Value = 0.0001396069650378791
$ws.Cells.Item(4, 16).Value = 0.0001396069650378791
$ws.Cells.Item(4, 17).Value = 0.1423664042
$ws.Cells.Item(4, 18).Value = 1.2812976378
$ws.Cells.Item(4, 19).Value = 0.00001462083667643619
$ws.Cells.Item(4, 20).Value = 0.00001462083667643619
$ws.Cells.Item(5, 7).Value = 17.41485066666667
$ws.Cells.Item(5, 8).Value = 52.24455200000001
$ws.Cells.Item(5, 9).Value = 0.1047285618770465
$ws.Cells.Item(5, 10).Value = 0.1047285618770465
$ws.Cells.Item(5, 13).Value = 58.099467
$ws.Cells.Item(5, 14).Value = 174.298401
$ws.Cells.Item(5, 15).Value = 0.992182294579622
$ws.Cells.Item(5, 16).Value = 0.992182294579622
$ws.Cells.Item(5, 17).Value = 1011.793541617928
$ws.Cells.Item(5, 18).Value = 9106.141874561354
$ws.Cells.Item(5, 19).Value = 0.103909824831192
$ws.Cells.Item(5, 20).Value = 0.103909824831192
$ws.Cells.Item(6, 9).Value = 0.1785014126970782
$ws.Cells.Item(6, 10).Value = 0.1785014126970782
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.06447966666666667
$ws.Cells.Item(6, 14).Value = 0.193439
$ws.Cells.Item(6, 15).Value = 0.001101138907643723
$ws.Cells.Item(6, 16).Value = 0.001101138907643722
$ws.Cells.Item(6, 17).Value = 1.913899178675778
$ws.Cells.Item(6, 18).Value = 17.225092608082
$ws.Cells.Item(6, 19).Value = 0.000196554850590122
$ws.Cells.Item(6, 20).Value = 0.000196554850590122
$ws.Cells.Item(7, 9).Value = 0.1785014126970782
$ws.Cells.Item(7, 10).Value = 0.1785014126970782
$ws.Cells.Item(7, 15).Value = 0.00657695954769643
$ws.Cells.Item(7, 16).Value = 0.006576959547696431
$ws.Cells.Item(7, 19).Value = 0.001173996570515349
$ws.Cells.Item(7, 20).Value = 0.001173996570515349
$ws.Cells.Item(8, 9).Value = 0.1785014126970782
$ws.Cells.Item(8, 10).Value = 0.1785014126970782
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.008175
$ws.Cells.Item(8, 14).Value = 0.024525
$ws.Cells.Item(8, 15).Value = 0.0001396069650378791
$ws.Cells.Item(8, 16).Value = 0.0001396069650378791
$ws.Cells.Item(8, 17).Value = 0.24265208855
$ws.Cells.Item(8, 18).Value = 2.18386879695
$ws.Cells.Item(8, 19).Value = 0.00002492004048161302
$ws.Cells.Item(8, 20).Value = 0.00002492004048161302
$ws.Cells.Item(9, 9).Value = 0.1785014126970782
$ws.Cells.Item(9, 10).Value = 0.1785014126970782
$ws.Cells.Item(9, 13).Value = 58.099467
$ws.Cells.Item(9, 14).Value = 174.298401
$ws.Cells.Item(9, 15).Value = 0.992182294579622
$ws.Cells.Item(9, 16).Value = 0.992182294579622
$ws.Cells.Item(9, 17).Value = 1724.520735313982
$ws.Cells.Item(9, 18).Value = 15520.68661782584
$ws.Cells.Item(9, 19).Value = 0.1771059412354911
$ws.Cells.Item(9, 20).Value = 0.1771059412354911
$ws.Cells.Item(10, 7).Value = 84.03051233333333
$ws.Cells.Item(10, 8).Value = 252.091537
$ws.Cells.Item(10, 9).Value = 0.5053385113032314
$ws.Cells.Item(10, 10).Value = 0.5053385113032314
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.06447966666666667
$ws.Cells.Item(10, 14).Value = 0.193439
$ws.Cells.Item(10, 15).Value = 0.001101138907643723
$ws.Cells.Item(10, 16).Value = 0.001101138907643722
$ws.Cells.Item(10, 17).Value = 5.418259425082556
$ws.Cells.Item(10, 18).Value = 48.764334825743
$ws.Cells.Item(10, 19).Value = 0.0005564478963267451
$ws.Cells.Item(10, 20).Value = 0.000556447896326745
$ws.Cells.Item(11, 7).Value = 84.03051233333333
$ws.Cells.Item(11, 8).Value = 252.091537
$ws.Cells.Item(11, 9).Value = 0.5053385113032314
$ws.Cells.Item(11, 10).Value = 0.5053385113032314
$ws.Cells.Item(11, 15).Value = 0.00657695954769643
$ws.Cells.Item(11, 16).Value = 0.006576959547696431
$ws.Cells.Item(11, 17).Value = 32.36255917425355
$ws.Cells.Item(11, 18).Value = 291.263032568282
$ws.Cells.Item(11, 19).Value = 0.003323590946734488
$ws.Cells.Item(11, 20).Value = 0.003323590946734489
$ws.Cells.Item(12, 7).Value = 84.03051233333333
$ws.Cells.Item(12, 8).Value = 252.091537
$ws.Cells.Item(12, 9).Value = 0.5053385113032314
$ws.Cells.Item(12, 10).Value = 0.5053385113032314
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.008175
$ws.Cells.Item(12, 14).Value = 0.024525
$ws.Cells.Item(12, 15).Value = 0.0001396069650378791
$ws.Cells.Item(12, 16).Value = 0.0001396069650378791
$ws.Cells.Item(12, 17).Value = 0.686949438325
$ws.Cells.Item(12, 18).Value = 6.182544944925001
$ws.Cells.Item(12, 19).Value = 0.00007054877587980408
$ws.Cells.Item(12, 20).Value = 0.00007054877587980408
$ws.Cells.Item(13, 7).Value = 84.03051233333333
$ws.Cells.Item(13, 8).Value = 252.091537
$ws.Cells.Item(13, 9).Value = 0.5053385113032314
$ws.Cells.Item(13, 10).Value = 0.5053385113032314
$ws.Cells.Item(13, 13).Value = 58.099467
$ws.Cells.Item(13, 14).Value = 174.298401
$ws.Cells.Item(13, 15).Value = 0.992182294579622
$ws.Cells.Item(13, 16).Value = 0.992182294579622
$ws.Cells.Item(13, 17).Value = 4882.127978303593
$ws.Cells.Item(13, 18).Value = 43939.15180473234
$ws.Cells.Item(13, 19).Value = 0.5013879236842904
$ws.Cells.Item(13, 20).Value = 0.5013879236842904
$ws.Cells.Item(14, 7).Value = 35.158014
$ws.Cells.Item(14, 8).Value = 105.474042
$ws.Cells.Item(14, 9).Value = 0.2114315141226439
$ws.Cells.Item(14, 10).Value = 0.2114315141226439
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.06447966666666667
$ws.Cells.Item(14, 14).Value = 0.193439
$ws.Cells.Item(14, 15).Value = 0.001101138907643723
$ws.Cells.Item(14, 16).Value = 0.001101138907643722
$ws.Cells.Item(14, 17).Value = 2.266977023382
$ws.Cells.Item(14, 18).Value = 20.402793210438
$ws.Cells.Item(14, 19).Value = 0.0002328154665024664
$ws.Cells.Item(14, 20).Value = 0.0002328154665024663
$ws.Cells.Item(15, 7).Value = 35.158014
$ws.Cells.Item(15, 8).Value = 105.474042
$ws.Cells.Item(15, 9).Value = 0.2114315141226439
$ws.Cells.Item(15, 10).Value = 0.2114315141226439
$ws.Cells.Item(15, 15).Value = 0.00657695954769643
$ws.Cells.Item(15, 16).Value = 0.006576959547696431
$ws.Cells.Item(15, 17).Value = 13.540359054468
$ws.Cells.Item(15, 18).Value = 121.863231490212
$ws.Cells.Item(15, 19).Value = 0.001390576515492835
$ws.Cells.Item(15, 20).Value = 0.001390576515492835
$ws.Cells.Item(16, 7).Value = 35.158014
$ws.Cells.Item(16, 8).Value = 105.474042
$ws.Cells.Item(16, 9).Value = 0.2114315141226439
$ws.Cells.Item(16, 10).Value = 0.2114315141226439
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.008175
$ws.Cells.Item(16, 14).Value = 0.024525
$ws.Cells.Item(16, 15).Value = 0.0001396069650378791
$ws.Cells.Item(16, 16).Value = 0.0001396069650378791
$ws.Cells.Item(16, 17).Value = 0.28741676445
$ws.Cells.Item(16, 18).Value = 2.58675088005
$ws.Cells.Item(16, 19).Value = 0.00002951731200002578
$ws.Cells.Item(16, 20).Value = 0.00002951731200002578
$ws.Cells.Item(17, 7).Value = 35.158014
$ws.Cells.Item(17, 8).Value = 105.474042
$ws.Cells.Item(17, 9).Value = 0.2114315141226439
$ws.Cells.Item(17, 10).Value = 0.2114315141226439
$ws.Cells.Item(17, 13).Value = 58.099467
$ws.Cells.Item(17, 14).Value = 174.298401
$ws.Cells.Item(17, 15).Value = 0.992182294579622
$ws.Cells.Item(17, 16).Value = 0.992182294579622
$ws.Cells.Item(17, 17).Value = 2042.661874178538
$ws.Cells.Item(17, 18).Value = 18383.95686760684
$ws.Cells.Item(17, 19).Value = 0.2097786048286486
$ws.Cells.Item(17, 20).Value = 0.2097786048286486
